# This script reproduces the edit described in the diff for
# Jogos_do_Dia_Betfair_Back_Lay_2026-01-02.xlsx:
#  - A new fixture row ("Algerian Ligue 1", Kabylie vs MC Alger) is inserted as
#    row 4, shifting the former rows 4-8 down to rows 5-9.
#  - Odds values across several rows are refreshed to updated quotes.
#  - The sheet dimension grows from A1:AO8 to A1:AO9.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 4 for the Algerian Ligue 1 match; this pushes the
# existing rows 4-8 down to rows 5-9.
$ws.Rows.Item(4).Insert()

# Row 2: Australian A-League Men | Melbourne Victory vs Perth Glory
$ws.Cells.Item(2, 1).Value = "Australian A-League Men"
$ws.Cells.Item(2, 2).Value = "'2026-01-02"
$ws.Cells.Item(2, 2).ClearFormats()
$ws.Cells.Item(2, 3).Value = "'05:35:00"
$ws.Cells.Item(2, 3).ClearFormats()
$ws.Cells.Item(2, 4).Value = "Melbourne Victory"
$ws.Cells.Item(2, 5).Value = "Perth Glory"
$ws.Cells.Item(2, 6).Value = 1.85
$ws.Cells.Item(2, 7).Value = 1.88
$ws.Cells.Item(2, 8).Value = 4.3
$ws.Cells.Item(2, 9).Value = 4.6
$ws.Cells.Item(2, 10).Value = 4.1
$ws.Cells.Item(2, 11).Value = 4.3
$ws.Cells.Item(2, 12).Value = 1.34
$ws.Cells.Item(2, 13).Value = 1.05
$ws.Cells.Item(2, 14).Value = 4.5
$ws.Cells.Item(2, 15).Value = 1.25
$ws.Cells.Item(2, 16).Value = 2.24
$ws.Cells.Item(2, 17).Value = 1.73
$ws.Cells.Item(2, 18).Value = 1.49
$ws.Cells.Item(2, 19).Value = 2.86
$ws.Cells.Item(2, 20).Value = 1.7
$ws.Cells.Item(2, 21).Value = 2.28
$ws.Cells.Item(2, 22).Value = 1.27
$ws.Cells.Item(2, 23).Value = 2.12
$ws.Cells.Item(2, 24).Value = 19.5
$ws.Cells.Item(2, 25).Value = 1000
$ws.Cells.Item(2, 26).Value = 40
$ws.Cells.Item(2, 27).Value = 100
$ws.Cells.Item(2, 28).Value = 11
$ws.Cells.Item(2, 29).Value = 9.4
$ws.Cells.Item(2, 30).Value = 18
$ws.Cells.Item(2, 31).Value = 55
$ws.Cells.Item(2, 32).Value = 12.5
$ws.Cells.Item(2, 33).Value = 9.6
$ws.Cells.Item(2, 34).Value = 17.5
$ws.Cells.Item(2, 35).Value = 55
$ws.Cells.Item(2, 36).Value = 21
$ws.Cells.Item(2, 37).Value = 18
$ws.Cells.Item(2, 38).Value = 30
$ws.Cells.Item(2, 39).Value = 95
$ws.Cells.Item(2, 40).Value = 10.5
$ws.Cells.Item(2, 41).Value = 46

# Row 3: Saudi Professional League | Al Najma Club vs Al-Khaleej Saihat
$ws.Cells.Item(3, 1).Value = "Saudi Professional League"
$ws.Cells.Item(3, 2).Value = "'2026-01-02"
$ws.Cells.Item(3, 2).ClearFormats()
$ws.Cells.Item(3, 3).Value = "'10:30:00"
$ws.Cells.Item(3, 3).ClearFormats()
$ws.Cells.Item(3, 4).Value = "Al Najma Club"
$ws.Cells.Item(3, 5).Value = "Al-Khaleej Saihat"
$ws.Cells.Item(3, 6).Value = 3.8
$ws.Cells.Item(3, 7).Value = 5.4
$ws.Cells.Item(3, 8).Value = 1.34
$ws.Cells.Item(3, 9).Value = 2.16
$ws.Cells.Item(3, 10).Value = 3.7
$ws.Cells.Item(3, 11).Value = 1000
$ws.Cells.Item(3, 12).Value = 1.01
$ws.Cells.Item(3, 13).Value = 1.01
$ws.Cells.Item(3, 14).Value = 2.02
$ws.Cells.Item(3, 15).Value = 1.22
$ws.Cells.Item(3, 16).Value = 2
$ws.Cells.Item(3, 17).Value = 1.64
$ws.Cells.Item(3, 18).Value = 1.38
$ws.Cells.Item(3, 19).Value = 2.42
$ws.Cells.Item(3, 20).Value = 1.01
$ws.Cells.Item(3, 21).Value = 1.01
$ws.Cells.Item(3, 22).Value = 1.87
$ws.Cells.Item(3, 23).Value = 1.23
$ws.Cells.Item(3, 24).Value = 1000
$ws.Cells.Item(3, 25).Value = 1000
$ws.Cells.Item(3, 26).Value = 1000
$ws.Cells.Item(3, 27).Value = 1000
$ws.Cells.Item(3, 28).Value = 1000
$ws.Cells.Item(3, 29).Value = 1000
$ws.Cells.Item(3, 30).Value = 1000
$ws.Cells.Item(3, 31).Value = 1000
$ws.Cells.Item(3, 32).Value = 1000
$ws.Cells.Item(3, 33).Value = 1000
$ws.Cells.Item(3, 34).Value = 1000
$ws.Cells.Item(3, 35).Value = 1000
$ws.Cells.Item(3, 36).Value = 1000
$ws.Cells.Item(3, 37).Value = 1000
$ws.Cells.Item(3, 38).Value = 1000
$ws.Cells.Item(3, 39).Value = 1000
$ws.Cells.Item(3, 40).Value = 1000
$ws.Cells.Item(3, 41).Value = 1000

# Row 4: Algerian Ligue 1 | Kabylie vs MC Alger
$ws.Cells.Item(4, 1).Value = "Algerian Ligue 1"
$ws.Cells.Item(4, 2).Value = "'2026-01-02"
$ws.Cells.Item(4, 2).ClearFormats()
$ws.Cells.Item(4, 3).Value = "'11:00:00"
$ws.Cells.Item(4, 3).ClearFormats()
$ws.Cells.Item(4, 4).Value = "Kabylie"
$ws.Cells.Item(4, 5).Value = "MC Alger"
$ws.Cells.Item(4, 6).Value = 1.04
$ws.Cells.Item(4, 7).Value = 1000
$ws.Cells.Item(4, 8).Value = 1.04
$ws.Cells.Item(4, 9).Value = 1000
$ws.Cells.Item(4, 10).Value = 1.01
$ws.Cells.Item(4, 11).Value = 1000
$ws.Cells.Item(4, 12).Value = 0
$ws.Cells.Item(4, 13).Value = 0
$ws.Cells.Item(4, 14).Value = 0
$ws.Cells.Item(4, 15).Value = 0
$ws.Cells.Item(4, 16).Value = 1.24
$ws.Cells.Item(4, 17).Value = 1.01
$ws.Cells.Item(4, 18).Value = 0
$ws.Cells.Item(4, 19).Value = 0
$ws.Cells.Item(4, 20).Value = 0
$ws.Cells.Item(4, 21).Value = 0
$ws.Cells.Item(4, 22).Value = 0
$ws.Cells.Item(4, 23).Value = 0
$ws.Cells.Item(4, 24).Value = 0
$ws.Cells.Item(4, 25).Value = 0
$ws.Cells.Item(4, 26).Value = 0
$ws.Cells.Item(4, 27).Value = 0
$ws.Cells.Item(4, 28).Value = 0
$ws.Cells.Item(4, 29).Value = 0
$ws.Cells.Item(4, 30).Value = 0
$ws.Cells.Item(4, 31).Value = 0
$ws.Cells.Item(4, 32).Value = 0
$ws.Cells.Item(4, 33).Value = 0
$ws.Cells.Item(4, 34).Value = 0
$ws.Cells.Item(4, 35).Value = 0
$ws.Cells.Item(4, 36).Value = 0
$ws.Cells.Item(4, 37).Value = 0
$ws.Cells.Item(4, 38).Value = 0
$ws.Cells.Item(4, 39).Value = 0
$ws.Cells.Item(4, 40).Value = 0
$ws.Cells.Item(4, 41).Value = 0

# Row 5: Saudi Professional League | Al-Ettifaq vs Al-Akhdoud
$ws.Cells.Item(5, 1).Value = "Saudi Professional League"
$ws.Cells.Item(5, 2).Value = "'2026-01-02"
$ws.Cells.Item(5, 2).ClearFormats()
$ws.Cells.Item(5, 3).Value = "'11:35:00"
$ws.Cells.Item(5, 3).ClearFormats()
$ws.Cells.Item(5, 4).Value = "Al-Ettifaq"
$ws.Cells.Item(5, 5).Value = "Al-Akhdoud"
$ws.Cells.Item(5, 6).Value = 1.76
$ws.Cells.Item(5, 7).Value = 2
$ws.Cells.Item(5, 8).Value = 4.4
$ws.Cells.Item(5, 9).Value = 6
$ws.Cells.Item(5, 10).Value = 3.7
$ws.Cells.Item(5, 11).Value = 4.4
$ws.Cells.Item(5, 12).Value = 0
$ws.Cells.Item(5, 13).Value = 0
$ws.Cells.Item(5, 14).Value = 0
$ws.Cells.Item(5, 15).Value = 0
$ws.Cells.Item(5, 16).Value = 2.02
$ws.Cells.Item(5, 17).Value = 1.78
$ws.Cells.Item(5, 18).Value = 0
$ws.Cells.Item(5, 19).Value = 0
$ws.Cells.Item(5, 20).Value = 0
$ws.Cells.Item(5, 21).Value = 0
$ws.Cells.Item(5, 22).Value = 0
$ws.Cells.Item(5, 23).Value = 0
$ws.Cells.Item(5, 24).Value = 0
$ws.Cells.Item(5, 25).Value = 0
$ws.Cells.Item(5, 26).Value = 0
$ws.Cells.Item(5, 27).Value = 0
$ws.Cells.Item(5, 28).Value = 0
$ws.Cells.Item(5, 29).Value = 0
$ws.Cells.Item(5, 30).Value = 0
$ws.Cells.Item(5, 31).Value = 0
$ws.Cells.Item(5, 32).Value = 0
$ws.Cells.Item(5, 33).Value = 0
$ws.Cells.Item(5, 34).Value = 0
$ws.Cells.Item(5, 35).Value = 0
$ws.Cells.Item(5, 36).Value = 0
$ws.Cells.Item(5, 37).Value = 0
$ws.Cells.Item(5, 38).Value = 0
$ws.Cells.Item(5, 39).Value = 0
$ws.Cells.Item(5, 40).Value = 0
$ws.Cells.Item(5, 41).Value = 0

# Row 6: Saudi Professional League | Al Ahli vs Al Nassr
$ws.Cells.Item(6, 1).Value = "Saudi Professional League"
$ws.Cells.Item(6, 2).Value = "'2026-01-02"
$ws.Cells.Item(6, 2).ClearFormats()
$ws.Cells.Item(6, 3).Value = "'14:30:00"
$ws.Cells.Item(6, 3).ClearFormats()
$ws.Cells.Item(6, 4).Value = "Al Ahli"
$ws.Cells.Item(6, 5).Value = "Al Nassr"
$ws.Cells.Item(6, 6).Value = 4.1
$ws.Cells.Item(6, 7).Value = 4.9
$ws.Cells.Item(6, 8).Value = 1.71
$ws.Cells.Item(6, 9).Value = 1.85
$ws.Cells.Item(6, 10).Value = 4.3
$ws.Cells.Item(6, 11).Value = 5.2
$ws.Cells.Item(6, 12).Value = 0
$ws.Cells.Item(6, 13).Value = 0
$ws.Cells.Item(6, 14).Value = 0
$ws.Cells.Item(6, 15).Value = 0
$ws.Cells.Item(6, 16).Value = 2.56
$ws.Cells.Item(6, 17).Value = 1.52
$ws.Cells.Item(6, 18).Value = 0
$ws.Cells.Item(6, 19).Value = 0
$ws.Cells.Item(6, 20).Value = 0
$ws.Cells.Item(6, 21).Value = 0
$ws.Cells.Item(6, 22).Value = 0
$ws.Cells.Item(6, 23).Value = 0
$ws.Cells.Item(6, 24).Value = 0
$ws.Cells.Item(6, 25).Value = 0
$ws.Cells.Item(6, 26).Value = 0
$ws.Cells.Item(6, 27).Value = 0
$ws.Cells.Item(6, 28).Value = 0
$ws.Cells.Item(6, 29).Value = 0
$ws.Cells.Item(6, 30).Value = 0
$ws.Cells.Item(6, 31).Value = 0
$ws.Cells.Item(6, 32).Value = 0
$ws.Cells.Item(6, 33).Value = 0
$ws.Cells.Item(6, 34).Value = 0
$ws.Cells.Item(6, 35).Value = 0
$ws.Cells.Item(6, 36).Value = 0
$ws.Cells.Item(6, 37).Value = 0
$ws.Cells.Item(6, 38).Value = 0
$ws.Cells.Item(6, 39).Value = 0
$ws.Cells.Item(6, 40).Value = 0
$ws.Cells.Item(6, 41).Value = 0

# Row 7: French Ligue 1 | Toulouse vs Lens
$ws.Cells.Item(7, 1).Value = "French Ligue 1"
$ws.Cells.Item(7, 2).Value = "'2026-01-02"
$ws.Cells.Item(7, 2).ClearFormats()
$ws.Cells.Item(7, 3).Value = "'16:45:00"
$ws.Cells.Item(7, 3).ClearFormats()
$ws.Cells.Item(7, 4).Value = "Toulouse"
$ws.Cells.Item(7, 5).Value = "Lens"
$ws.Cells.Item(7, 6).Value = 2.9
$ws.Cells.Item(7, 7).Value = 2.96
$ws.Cells.Item(7, 8).Value = 2.66
$ws.Cells.Item(7, 9).Value = 2.72
$ws.Cells.Item(7, 10).Value = 3.4
$ws.Cells.Item(7, 11).Value = 3.5
$ws.Cells.Item(7, 12).Value = 0
$ws.Cells.Item(7, 13).Value = 1.08
$ws.Cells.Item(7, 14).Value = 3.7
$ws.Cells.Item(7, 15).Value = 1.34
$ws.Cells.Item(7, 16).Value = 1.9
$ws.Cells.Item(7, 17).Value = 2
$ws.Cells.Item(7, 18).Value = 1.36
$ws.Cells.Item(7, 19).Value = 3.6
$ws.Cells.Item(7, 20).Value = 1.77
$ws.Cells.Item(7, 21).Value = 2.2
$ws.Cells.Item(7, 22).Value = 0
$ws.Cells.Item(7, 23).Value = 0
$ws.Cells.Item(7, 24).Value = 13.5
$ws.Cells.Item(7, 25).Value = 11
$ws.Cells.Item(7, 26).Value = 22
$ws.Cells.Item(7, 27).Value = 42
$ws.Cells.Item(7, 28).Value = 12
$ws.Cells.Item(7, 29).Value = 7.6
$ws.Cells.Item(7, 30).Value = 12.5
$ws.Cells.Item(7, 31).Value = 32
$ws.Cells.Item(7, 32).Value = 23
$ws.Cells.Item(7, 33).Value = 13
$ws.Cells.Item(7, 34).Value = 16.5
$ws.Cells.Item(7, 35).Value = 42
$ws.Cells.Item(7, 36).Value = 46
$ws.Cells.Item(7, 37).Value = 38
$ws.Cells.Item(7, 38).Value = 46
$ws.Cells.Item(7, 39).Value = 100
$ws.Cells.Item(7, 40).Value = 32
$ws.Cells.Item(7, 41).Value = 26

# Row 8: Italian Serie A | Cagliari vs AC Milan
$ws.Cells.Item(8, 1).Value = "Italian Serie A"
$ws.Cells.Item(8, 2).Value = "'2026-01-02"
$ws.Cells.Item(8, 2).ClearFormats()
$ws.Cells.Item(8, 3).Value = "'16:45:00"
$ws.Cells.Item(8, 3).ClearFormats()
$ws.Cells.Item(8, 4).Value = "Cagliari"
$ws.Cells.Item(8, 5).Value = "AC Milan"
$ws.Cells.Item(8, 6).Value = 7
$ws.Cells.Item(8, 7).Value = 7.4
$ws.Cells.Item(8, 8).Value = 1.56
$ws.Cells.Item(8, 9).Value = 1.58
$ws.Cells.Item(8, 10).Value = 4.4
$ws.Cells.Item(8, 11).Value = 4.5
$ws.Cells.Item(8, 12).Value = 0
$ws.Cells.Item(8, 13).Value = 1.07
$ws.Cells.Item(8, 14).Value = 3.85
$ws.Cells.Item(8, 15).Value = 1.32
$ws.Cells.Item(8, 16).Value = 1.97
$ws.Cells.Item(8, 17).Value = 1.97
$ws.Cells.Item(8, 18).Value = 1.37
$ws.Cells.Item(8, 19).Value = 3.5
$ws.Cells.Item(8, 20).Value = 2.08
$ws.Cells.Item(8, 21).Value = 1.87
$ws.Cells.Item(8, 22).Value = 0
$ws.Cells.Item(8, 23).Value = 0
$ws.Cells.Item(8, 24).Value = 16
$ws.Cells.Item(8, 25).Value = 7.8
$ws.Cells.Item(8, 26).Value = 8.6
$ws.Cells.Item(8, 27).Value = 14.5
$ws.Cells.Item(8, 28).Value = 22
$ws.Cells.Item(8, 29).Value = 9.800000000000001
$ws.Cells.Item(8, 30).Value = 10.5
$ws.Cells.Item(8, 31).Value = 17.5
$ws.Cells.Item(8, 32).Value = 65
$ws.Cells.Item(8, 33).Value = 28
$ws.Cells.Item(8, 34).Value = 27
$ws.Cells.Item(8, 35).Value = 42
$ws.Cells.Item(8, 36).Value = 280
$ws.Cells.Item(8, 37).Value = 1000
$ws.Cells.Item(8, 38).Value = 130
$ws.Cells.Item(8, 39).Value = 190
$ws.Cells.Item(8, 40).Value = 1000
$ws.Cells.Item(8, 41).Value = 9.199999999999999

# Row 9: Spanish La Liga | Rayo Vallecano vs Getafe
$ws.Cells.Item(9, 1).Value = "Spanish La Liga"
$ws.Cells.Item(9, 2).Value = "'2026-01-02"
$ws.Cells.Item(9, 2).ClearFormats()
$ws.Cells.Item(9, 3).Value = "'17:00:00"
$ws.Cells.Item(9, 3).ClearFormats()
$ws.Cells.Item(9, 4).Value = "Rayo Vallecano"
$ws.Cells.Item(9, 5).Value = "Getafe"
$ws.Cells.Item(9, 6).Value = 2.24
$ws.Cells.Item(9, 7).Value = 2.26
$ws.Cells.Item(9, 8).Value = 4.3
$ws.Cells.Item(9, 9).Value = 4.5
$ws.Cells.Item(9, 10).Value = 3
$ws.Cells.Item(9, 11).Value = 3.1
$ws.Cells.Item(9, 12).Value = 0
$ws.Cells.Item(9, 13).Value = 1.17
$ws.Cells.Item(9, 14).Value = 2.3
$ws.Cells.Item(9, 15).Value = 1.73
$ws.Cells.Item(9, 16).Value = 1.41
$ws.Cells.Item(9, 17).Value = 3.3
$ws.Cells.Item(9, 18).Value = 1.14
$ws.Cells.Item(9, 19).Value = 7.6
$ws.Cells.Item(9, 20).Value = 2.6
$ws.Cells.Item(9, 21).Value = 1.58
$ws.Cells.Item(9, 22).Value = 0
$ws.Cells.Item(9, 23).Value = 0
$ws.Cells.Item(9, 24).Value = 7
$ws.Cells.Item(9, 25).Value = 9.800000000000001
$ws.Cells.Item(9, 26).Value = 32
$ws.Cells.Item(9, 27).Value = 160
$ws.Cells.Item(9, 28).Value = 6
$ws.Cells.Item(9, 29).Value = 7.4
$ws.Cells.Item(9, 30).Value = 22
$ws.Cells.Item(9, 31).Value = 110
$ws.Cells.Item(9, 32).Value = 10.5
$ws.Cells.Item(9, 33).Value = 13.5
$ws.Cells.Item(9, 34).Value = 34
$ws.Cells.Item(9, 35).Value = 170
$ws.Cells.Item(9, 36).Value = 30
$ws.Cells.Item(9, 37).Value = 38
$ws.Cells.Item(9, 38).Value = 95
$ws.Cells.Item(9, 39).Value = 430
$ws.Cells.Item(9, 40).Value = 40
$ws.Cells.Item(9, 41).Value = 240
